# Update the EC (Estado de Cuenta) workbook:
# - Remove the first worker data row (CC 73202748 - UMERLIS ANTONIO MARQUEZ GUEVARA)
# - Refresh the summary totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) to
#   reflect the remaining two worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 16 (first data row), shifting the remaining rows up.
$ws.Rows("16").Delete()

# Recalculate / update summary values that depended on the deleted row.
$ws.Range("E11").Value = 3258   # Valor Mora total (1525 + 1733)
$ws.Range("C13").Value = 2      # Cant. Trabajadores (was 3)
$ws.Range("F13").Value = 1      # Cant. Periodos (was 2)
